# Apply the "Updated symbol list" price/coin refresh to the crypto table.
# Column D holds numeric-looking values that must remain stored as text
# (as in the source data), so NumberFormat is forced to "@" (Text) before
# assigning each such value; purely textual cells (B/C/E) don't need that.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "271.12"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.66"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.340"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06199"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.658"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.661"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.399"
$ws.Range("E8").Value = "7FTXTokenFTT"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8316"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01373"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1611"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08296"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03565"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03219"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.066"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.09300"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001640"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04751"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006379"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005678"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001502"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.730"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3334"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1235"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002708"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04741"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007040"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1162"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003566"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01189"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006274"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0009920"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.7833"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002324"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002404"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01242"
